$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.879.59"
$ws.Range("E2").Value = "  +4.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.778.55"
$ws.Range("E3").Value = "  +4.44%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.13"
$ws.Range("E5").Value = "  +4.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.33"
$ws.Range("E6").Value = "  +1.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  +4.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  +4.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.58"
$ws.Range("E10").Value = "  +6.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0853"
$ws.Range("E11").Value = "  +4.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.04"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.214.58"
$ws.Range("E15").Value = "  +4.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.777.39"
$ws.Range("E16").Value = "  +4.90%  "

$ws.Range("E17").Value = "  +1.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.809.32"
$ws.Range("E18").Value = "  +4.04%  "

$ws.Range("E19").Value = "  +10.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  +4.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("E21").Value = "  -3.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +2.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.20"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.97"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("E25").Value = "  +7.62%  "

$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.62"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("E33").Value = "  +3.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0819"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("E37").Value = "  +2.24%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0391"
$ws.Range("E39").Value = "  +12.07%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  +1.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +23.99%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  +3.10%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.93"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.34"
$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.065.56"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.54"
$ws.Range("E49").Value = "  +4.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.900"
$ws.Range("E50").Value = "  +13.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.86"
$ws.Range("E51").Value = "  -1.64%  "
